$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'youth black knee pads'
$ws.Range("A2").Value = 'basketball clothes for boys'
$ws.Range("A3").Value = 'men capri'
$ws.Range("A4").Value = 'black baseball pants adult small'
$ws.Range("A5").Value = 'under pants for men'
$ws.Range("A6").Value = 'by design knee pads'
$ws.Range("A7").Value = 'baseball knee high pants mens'
$ws.Range("A8").Value = 'boys youth compression tights'
$ws.Range("A9").Value = 'mens small running tights'
$ws.Range("A10").Value = 'youth large softball pants'
$ws.Range("A11").Value = 'small basketballs'
$ws.Range("A12").Value = 'cycling capri'
$ws.Range("A13").Value = 'mens black basketball shorts'
$ws.Range("A14").Value = 'knee protector for running'
$ws.Range("A15").Value = 'knees protection'
$ws.Range("A16").Value = 'dry fit baseball pants'
$ws.Range("A17").Value = 'girls black baseball pants'
$ws.Range("A18").Value = 'mens baseball sliding shorts'
$ws.Range("A19").Value = 'wrestling equipment'
$ws.Range("A20").Value = 'the knee pads'
$ws.Range("A21").Value = 'football compression'
$ws.Range("A22").Value = 'soccer guards for men'
$ws.Range("A23").Value = 'spandex leggings for boys'
$ws.Range("A24").Value = 'basketball pants men'
$ws.Range("A25").Value = 'hunting knee pads'
$ws.Range("A26").Value = 'combat pants knee pads'
$ws.Range("A27").Value = 'nike compression basketball pants'
$ws.Range("A28").Value = 'seamless capri leggings'
$ws.Range("A29").Value = 'basketball knee pads adidas'
$ws.Range("A30").Value = 'women black leggings'
$ws.Range("A31").Value = 'red basketball knee pads'
$ws.Range("A32").Value = 'underarmour mens leggings'
$ws.Range("A33").Value = 'red nike compression pants men'
$ws.Range("A34").Value = 'men compression pants long'
$ws.Range("A35").Value = 'knee tights'
$ws.Range("A36").Value = 'knee protector basketball'
$ws.Range("A37").Value = 'youth football pants with pads'
$ws.Range("A38").Value = 'little boys compression tights'
$ws.Range("A39").Value = 'wrestling tights for boys'
$ws.Range("A40").Value = 'youth boys baseball pants'
$ws.Range("A41").Value = 'padded basketball tights'
$ws.Range("A42").Value = 'basketball clothes youth'
$ws.Range("A43").Value = 'leggings boys'
$ws.Range("A44").Value = 'youth compression pants'
$ws.Range("A45").Value = 'knee pad youth'
$ws.Range("A46").Value = 'best knee pads for volleyball'
$ws.Range("A47").Value = 'athletic pants for men big and tall'
$ws.Range("A48").Value = 'leg compression basketball'
$ws.Range("A49").Value = 'compression shorts baseball'
$ws.Range("A50").Value = 'pants for men sports'
$ws.Range("A51").Value = 'knee pads large'
$ws.Range("A52").Value = 'gym tights for men'
$ws.Range("A53").Value = 'soccer gear for men'
$ws.Range("A54").Value = 'large tall athletic pants men'
$ws.Range("A55").Value = 'youth football pads'
$ws.Range("A56").Value = 'padded compression shorts'
$ws.Range("A57").Value = 'youth girls softball pants black'
$ws.Range("A58").Value = 'athletic pants for men'
$ws.Range("A59").Value = 'patella knee'
$ws.Range("A60").Value = 'youth softball pants'
$ws.Range("A61").Value = 'basketball knee sleeves'
$ws.Range("A62").Value = 'knee pads for crossfit'
$ws.Range("A63").Value = 'mens nike dri fit compression pants'
$ws.Range("A64").Value = 'volleyball youth knee pads'
$ws.Range("A65").Value = 'troll knee pads'
$ws.Range("A66").Value = 'elastic knee pads'
$ws.Range("A67").Value = 'men leggings adidas'
$ws.Range("A68").Value = 'snowmobile knee pads'
$ws.Range("A69").Value = 'knee pads military'
$ws.Range("A70").Value = 'knee pads sleeve'
$ws.Range("A71").Value = 'knee pads mma'
$ws.Range("A72").Value = 'under armour compression pants men'
$ws.Range("A73").Value = 'mcdavid knee pads basketball'
$ws.Range("A74").Value = 'black leggings xsmall'
$ws.Range("A75").Value = 'tesla compression pants'
$ws.Range("A76").Value = 'fox knee pads'
$ws.Range("A77").Value = 'biking knee pads'
$ws.Range("A78").Value = 'bmx knee pads'
$ws.Range("A79").Value = 'mcdavid knee pad'
$ws.Range("A80").Value = 'pant with knee pads'
$ws.Range("A81").Value = 'green knee pads'
$ws.Range("A82").Value = 'knee pads tights'
$ws.Range("A83").Value = 'warm compression pants'
$ws.Range("A84").Value = 'dancer knee pads'
$ws.Range("A85").Value = 'navy compression pants'
$ws.Range("A86").Value = 'mizuno knee pad'
$ws.Range("A87").Value = 'mens pants with knee pads'
$ws.Range("A88").Value = 'black athletic capri'
$ws.Range("A89").Value = 'tactical knee pad'
$ws.Range("A90").Value = 'men compression pants blue'
$ws.Range("A91").Value = 'leggins for men sport'
$ws.Range("A92").Value = 'knee pads ski'
$ws.Range("A93").Value = 'knee pad snowboard'
$ws.Range("A94").Value = 'knee pads mcdavid basketball'
$ws.Range("A95").Value = 'pantalones con rodilleras'
$ws.Range("A96").Value = 'pantalon con rodilleras'
$ws.Range("A97").Value = 'cat knee pad pants'
$ws.Range("A98").Value = 'knee padded pants men'
$ws.Range("A99").Value = 'knee pad pants men'
$ws.Range("A100").Value = 'caterpillar knee pad pants'
